$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-02 Wednesday" "2024-10-03 Thursday"

Replace-Text "214÷2=107, 0" "772÷7=110, 2"
Replace-Text "161÷2=80, 1" "775÷2=387, 1"
Replace-Text "705÷5=141, 0" "915÷4=228, 3"
Replace-Text "456÷5=91, 1" "737÷9=81, 8"
Replace-Text "384÷8=48, 0" "761÷4=190, 1"

Replace-Text "544÷5=108, 4" "963÷2=481, 1"
Replace-Text "112÷8=14, 0" "431÷8=53, 7"
Replace-Text "627÷7=89, 4" "183÷9=20, 3"
Replace-Text "391÷3=130, 1" "503÷8=62, 7"
Replace-Text "878÷9=97, 5" "623÷5=124, 3"

Replace-Text "509÷8=63, 5" "925÷8=115, 5"
Replace-Text "973÷7=139, 0" "680÷5=136, 0"
Replace-Text "250÷8=31, 2" "131÷9=14, 5"
Replace-Text "465÷4=116, 1" "934÷5=186, 4"
Replace-Text "174÷9=19, 3" "499÷5=99, 4"

Replace-Text "425÷4=106, 1" "348÷6=58, 0"
Replace-Text "721÷2=360, 1" "717÷3=239, 0"
Replace-Text "539÷6=89, 5" "482÷3=160, 2"
Replace-Text "623÷4=155, 3" "584÷5=116, 4"
Replace-Text "409÷4=102, 1" "141÷2=70, 1"

Replace-Text "669÷3=223, 0" "599÷2=299, 1"
Replace-Text "118÷9=13, 1" "389÷3=129, 2"
Replace-Text "470÷7=67, 1" "987÷9=109, 6"
Replace-Text "868÷6=144, 4" "729÷2=364, 1"
Replace-Text "887÷4=221, 3" "187÷8=23, 3"

Write-Output "Done"
